# Weekly driver report update for 2025-05-05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" summary block (rows 3-5) ---
$ws.Range("D3").Value = 97.1

$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 621
$ws.Range("D4").Value = 97.6

$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 622

# --- "Good Drivers" detail block (rows 13-35) ---
# Each entry only lists the columns that actually change for that row.
$rows = @(
    @{ Row = 13; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.0.4"; B = 1869842; C = 3420; D = 4386; E = 2436; F = 1877648; H = "22.250.0.4"; J = "2023-07-25" },
    @{ Row = 14; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.40.0.7"; B = 8170878; C = 15867; D = 13188; E = 18436; F = 8199933; H = "22.40.0.7"; J = "2021-09-18" },
    @{ Row = 15; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.190.0.4"; B = 1611822; C = 4793; D = 1136; E = 2996; F = 1617751; H = "22.190.0.4"; J = "2022-11-22" },
    @{ Row = 16; B = 298304; C = 925; E = 420; F = 299532 },
    @{ Row = 17; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.10.0.8"; B = 381616; C = 1106; D = 582; E = 601; F = 383304; H = "23.10.0.8"; I = 99.6; J = "2023-10-30" },
    @{ Row = 18; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.90.0.5"; B = 335610; C = 495; D = 530; E = 414; F = 336635; H = "22.90.0.5"; J = "2021-09-26" },
    @{ Row = 19; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.220.0.4"; B = 750778; C = 1655; D = 507; E = 1290; F = 752940; H = "22.220.0.4"; J = "2023-03-28" },
    @{ Row = 20; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.20.0.6"; B = 1021828; C = 2379; D = 1006; E = 1021; F = 1025213; H = "22.20.0.6"; J = "2020-11-29" },
    @{ Row = 21; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.60.0.6"; B = 100380; C = 264; D = 30; E = 267; F = 100674; H = "22.60.0.6"; J = "2021-05-26" },
    @{ Row = 22; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.1.3"; B = 195464; C = 430; D = 70; E = 401; F = 195964; H = "23.80.1.3"; J = "2024-09-03" },
    @{ Row = 23; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.0.0.6"; B = 1833336; C = 2376; D = 656; E = 1899; F = 1836368; G = "intel(r) wi-fi 6 ax201 160mhz"; H = "22.0.0.6"; J = "2020-09-16" },
    @{ Row = 24; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.170.0.3"; B = 1237156; C = 2499; D = 583; E = 1951; F = 1240238; H = "22.170.0.3"; J = "2022-08-28" },
    @{ Row = 25; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.80.0.4"; B = 63994; C = 66; D = 41; E = 68; F = 64101; H = "21.80.0.4"; J = "2020-01-29" },
    @{ Row = 26; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"; B = 81417; C = 107; D = 25; E = 158; F = 81549; H = "22.80.0.9"; J = "2021-08-18" },
    @{ Row = 27; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.0.1.1"; B = 1400496; C = 2029; D = 1075; E = 3210; F = 1403600; H = "22.0.1.1"; J = "2020-09-28" },
    @{ Row = 28; A = "Marvell AVASTAR Wireless-AC Network Controller - 15.68.17022.122"; B = 293859; C = 502; D = 195; E = 668; F = 294556; G = "marvell avastar wireless-ac network controller"; H = "15.68.17022.122"; J = "2021-09-08" },
    @{ Row = 29; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.80.2.1"; B = 393045; C = 453; D = 216; E = 1251; F = 393714; H = "21.80.2.1"; I = 99.8; J = "2020-02-24" },
    @{ Row = 30; B = 509717; C = 522; E = 590; F = 510417 },
    @{ Row = 31; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"; B = 11362; C = 7; D = 0; E = 30; F = 11369; H = "21.40.1.3"; I = 99.9 },
    @{ Row = 32; A = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"; B = 75454; C = 2; D = 1; E = 79; F = 75457; H = "21.110.3.2"; J = "2020-08-05" },
    @{ Row = 34; B = 121310; C = 26; E = 154; F = 121350 },
    @{ Row = 35; B = 35363; E = 81; F = 35377 }
)

foreach ($r in $rows) {
    if ($r.ContainsKey("A")) { $ws.Cells.Item($r.Row, 1).Value = $r.A }
    if ($r.ContainsKey("B")) { $ws.Cells.Item($r.Row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { $ws.Cells.Item($r.Row, 4).Value = $r.D }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($r.Row, 5).Value = $r.E }
    if ($r.ContainsKey("F")) { $ws.Cells.Item($r.Row, 6).Value = $r.F }
    if ($r.ContainsKey("G")) { $ws.Cells.Item($r.Row, 7).Value = $r.G }
    if ($r.ContainsKey("H")) { $ws.Cells.Item($r.Row, 8).Value = $r.H }
    if ($r.ContainsKey("I")) { $ws.Cells.Item($r.Row, 9).Value = $r.I }
    if ($r.ContainsKey("J")) { $ws.Cells.Item($r.Row, 10).Value = "'" + $r.J }
}

# Row 31's driver-vintage date is cleared entirely (swapped with row 32's date above).
$ws.Range("J31").Value = ""
